# Update cryptocurrency price (column D) and volume change (column E) values
# to reflect the latest scrape performed by the GitHub Actions workflow.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'334.89"
$ws.Range("E2").Value = "'1.59%"
$ws.Range("D3").Value = "'43.98"
$ws.Range("E3").Value = "'6.08%"
$ws.Range("D4").Value = "'5.751"
$ws.Range("E4").Value = "'2.00%"
$ws.Range("D5").Value = "'0.08382"
$ws.Range("E5").Value = "'0.80%"
$ws.Range("D6").Value = "'8.868"
$ws.Range("E6").Value = "'1.13%"
$ws.Range("D7").Value = "'1.953"
$ws.Range("E7").Value = "'-4.36%"
$ws.Range("E8").Value = "'-4.12%"
$ws.Range("D9").Value = "'0.9480"
$ws.Range("E9").Value = "'2.54%"
$ws.Range("E10").Value = "'-2.96%"
$ws.Range("D11").Value = "'0.1987"
$ws.Range("E11").Value = "'1.46%"
$ws.Range("D12").Value = "'0.1004"
$ws.Range("E12").Value = "'6.79%"
$ws.Range("D13").Value = "'0.04404"
$ws.Range("E13").Value = "'11.10%"
$ws.Range("D14").Value = "'0.1069"
$ws.Range("E14").Value = "'0.76%"
$ws.Range("D15").Value = "'0.001298"
$ws.Range("E15").Value = "'-0.59%"
$ws.Range("D16").Value = "'0.006059"
$ws.Range("E16").Value = "'-1.46%"
$ws.Range("D17").Value = "'3.486"
$ws.Range("E17").Value = "'1.21%"
$ws.Range("D18").Value = "'4.525"
$ws.Range("E18").Value = "'-0.07%"
$ws.Range("E19").Value = "'0.19%"
$ws.Range("D20").Value = "'8.695"
$ws.Range("E20").Value = "'3.71%"
$ws.Range("D21").Value = "'0.1363"
$ws.Range("E21").Value = "'-0.78%"
$ws.Range("D23").Value = "'0.04405"
$ws.Range("E23").Value = "'0.23%"
$ws.Range("E24").Value = "'0.05%"
$ws.Range("D25").Value = "'0.004349"
$ws.Range("E25").Value = "'0.74%"
$ws.Range("D26").Value = "'0.0001264"
$ws.Range("E26").Value = "'5.33%"
$ws.Range("D27").Value = "'0.0004002"
$ws.Range("D39").Value = "'0.02811"
$ws.Range("E39").Value = "'1.11%"
$ws.Range("D40").Value = "'0.05914"
$ws.Range("E40").Value = "'7.73%"
$ws.Range("D41").Value = "'0.007960"
$ws.Range("E41").Value = "'0.69%"
$ws.Range("E42").Value = "'0.51%"
$ws.Range("D43").Value = "'0.009053"
$ws.Range("E43").Value = "'1.34%"
$ws.Range("D44").Value = "'0.002149"
$ws.Range("E44").Value = "'0.39%"
$ws.Range("D45").Value = "'0.01036"
$ws.Range("E45").Value = "'-12.12%"
$ws.Range("D46").Value = "'0.00007239"
$ws.Range("E46").Value = "'3.47%"
$ws.Range("E47").Value = "'0.30%"
$ws.Range("D48").Value = "'0.003205"
$ws.Range("E48").Value = "'0.43%"
$ws.Range("D49").Value = "'0.002277"
$ws.Range("E49").Value = "'-0.13%"
$ws.Range("D50").Value = "'0.00002107"
$ws.Range("E50").Value = "'0.30%"
$ws.Range("D51").Value = "'0.0002007"
$ws.Range("E51").Value = "'0.30%"
